# Replace Product, IT, and Finance templates with correct industry-specific
# content: this workbook is the "Product Development" comprehensive budget
# template, but several cells still carried over "Artificial Intelligence
# and Machine Learning" / "AI" / "ML" wording from another template it was
# cloned from. Swap that copy out for the correct Product Development text.

$wb = $excel.ActiveWorkbook

# --- Sheet: Instructions & User Guide ---------------------------------
$ws = $wb.Worksheets.Item("Instructions & User Guide")
$ws.Range("A1").Value = "Product Development Comprehensive Budget - User Guide & Instructions"
$ws.Range("A56").Value = "📋 PRODUCT DEVELOPMENT PROJECT OVERVIEW"
$ws.Range("B59").Value = "Data Scientists, Product Engineers, Product Architects, DevOps Engineers..."

# --- Sheet: Budget Summary ---------------------------------------------
$ws = $wb.Worksheets.Item("Budget Summary")
$ws.Range("A1").Value = "Product Development - Executive Budget Summary"

# --- Sheet: Resources ----------------------------------------------------
$ws = $wb.Worksheets.Item("Resources")
$ws.Range("A1").Value = "Product Development - Resources Budget"
$ws.Range("A5").Value = "Product Engineers"
$ws.Range("A6").Value = "Product Architects"

# --- Sheet: Logistics ----------------------------------------------------
$ws = $wb.Worksheets.Item("Logistics")
$ws.Range("A1").Value = "Product Development - Logistics Budget"

# --- Sheet: Technology ---------------------------------------------------
$ws = $wb.Worksheets.Item("Technology")
$ws.Range("A1").Value = "Product Development - Technology Budget"
$ws.Range("A5").Value = "Product Platform Licenses"

# --- Sheet: Training ------------------------------------------------------
$ws = $wb.Worksheets.Item("Training")
$ws.Range("A1").Value = "Product Development - Training Budget"
$ws.Range("A4").Value = "Product Development Certification Programs"
$ws.Range("A10").Value = "TOTAL TRProductNING"

# --- Sheet: Contingency ----------------------------------------------------
$ws = $wb.Worksheets.Item("Contingency")
$ws.Range("A1").Value = "Product Development - Contingency Budget"

# --- Sheet: Timeline ------------------------------------------------------
$ws = $wb.Worksheets.Item("Timeline")
$ws.Range("A1").Value = "Product Development - Budget Timeline"
